$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 72522.14
$ws.Range("J17").Value = 72522.14
$ws.Range("L17").Value = 217566.42
$ws.Range("N17").Value = -217902.42
$ws.Range("H19").Value = 756.3
$ws.Range("I19").Value = 746.5
$ws.Range("J19").Value = 771
$ws.Range("K19").Value = 746.5
$ws.Range("L19").Value = 771
$ws.Range("M19").Value = -571.5
$ws.Range("N19").Value = -1121
$ws.Range("H28").Value = 1554.7778
$ws.Range("I28").Value = 2297.4
$ws.Range("J28").Value = 626.5
$ws.Range("K28").Value = 2297.4
$ws.Range("L28").Value = 626.5
$ws.Range("M28").Value = -1812.4
$ws.Range("N28").Value = -1596.5
$ws.Range("H32").Value = 3886.111
$ws.Range("J32").Value = 2596
$ws.Range("L32").Value = 2596
$ws.Range("N32").Value = -3248
$ws.Range("H70").Value = 63749.688
$ws.Range("I70").Value = 201039
$ws.Range("J70").Value = 1345.4546
$ws.Range("K70").Value = 603117
$ws.Range("L70").Value = 4036.3638
$ws.Range("M70").Value = -602847
$ws.Range("N70").Value = -4576.3638
$ws.Range("H73").Value = 63749.688
$ws.Range("I73").Value = 201039
$ws.Range("J73").Value = 1345.4546
$ws.Range("K73").Value = 603117
$ws.Range("L73").Value = 4036.3638
$ws.Range("M73").Value = -602181
$ws.Range("N73").Value = -5908.3638
$ws.Range("H74").Value = 4506
$ws.Range("I74").Value = 4189.5557
$ws.Range("J74").Value = 4862
$ws.Range("K74").Value = 4189.5557
$ws.Range("L74").Value = 4862
$ws.Range("M74").Value = -3253.5557
$ws.Range("N74").Value = -6734
$ws.Range("H77").Value = 4506
$ws.Range("I77").Value = 4189.5557
$ws.Range("J77").Value = 4862
$ws.Range("K77").Value = 20947.7785
$ws.Range("L77").Value = 24310
$ws.Range("M77").Value = -16267.7785
$ws.Range("N77").Value = -33670
$ws.Range("H137").Value = 1556.8096
$ws.Range("I137").Value = 1411.875
$ws.Range("J137").Value = 2020.6
$ws.Range("K137").Value = 4235.625
$ws.Range("L137").Value = 6061.799999999999
$ws.Range("M137").Value = -1685.625
$ws.Range("N137").Value = -11161.8
$ws.Range("H138").Value = 2169.7578
$ws.Range("I138").Value = 1103.4897
$ws.Range("J138").Value = 3305.5652
$ws.Range("K138").Value = 3310.4691
$ws.Range("L138").Value = 9916.695599999999
$ws.Range("M138").Value = 1829.5309
$ws.Range("N138").Value = -20196.6956
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 33485.098
$ws.Range("I2").Value = 1035.3846
$ws.Range("J2").Value = 56921
$ws.Range("K2").Value = 1035.3846
$ws.Range("L2").Value = 56921
$ws.Range("M2").Value = -922.3846000000001
$ws.Range("N2").Value = -57147
$ws.Range("H9").Value = 33800
$ws.Range("J9").Value = 29750
$ws.Range("L9").Value = 29750
$ws.Range("N9").Value = -30090
$ws.Range("H20").Value = 33800
$ws.Range("J20").Value = 29750
$ws.Range("L20").Value = 29750
$ws.Range("N20").Value = -30290
$ws.Range("H23").Value = 11714.777
$ws.Range("I23").Value = 11006
$ws.Range("J23").Value = 11803.375
$ws.Range("K23").Value = 11006
$ws.Range("L23").Value = 11803.375
$ws.Range("M23").Value = -10747
$ws.Range("N23").Value = -12321.375
$ws.Range("H32").Value = 14152.046
$ws.Range("I32").Value = 14721.526
$ws.Range("J32").Value = 10545.333
$ws.Range("K32").Value = 14721.526
$ws.Range("L32").Value = 10545.333
$ws.Range("M32").Value = -14434.526
$ws.Range("N32").Value = -11119.333
$ws.Range("H63").Value = 4097
$ws.Range("I63").Value = 4097
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 4097
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -3411
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 4097
$ws.Range("I66").Value = 4097
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 20485
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -17053
$ws.Range("N66").ClearContents()
$ws.Range("H74").Value = 1296.8823
$ws.Range("I74").Value = 1153.9166
$ws.Range("J74").Value = 1640
$ws.Range("K74").Value = 1153.9166
$ws.Range("L74").Value = 1640
$ws.Range("M74").Value = -279.9166
$ws.Range("N74").Value = -3388
$ws.Range("H77").Value = 1296.8823
$ws.Range("I77").Value = 1153.9166
$ws.Range("J77").Value = 1640
$ws.Range("K77").Value = 5769.583000000001
$ws.Range("L77").Value = 8200
$ws.Range("M77").Value = -1401.583000000001
$ws.Range("N77").Value = -16936
$ws.Range("H110").Value = 1054.3334
$ws.Range("I110").Value = 1062.875
$ws.Range("J110").Value = 1044.5714
$ws.Range("K110").Value = 1062.875
$ws.Range("L110").Value = 1044.5714
$ws.Range("M110").Value = 982.125
$ws.Range("N110").Value = -5134.5714
$ws.Range("H116").Value = 33485.098
$ws.Range("I116").Value = 1035.3846
$ws.Range("J116").Value = 56921
$ws.Range("K116").Value = 1035.3846
$ws.Range("L116").Value = 56921
$ws.Range("M116").Value = 1258.6154
$ws.Range("N116").Value = -61509
$ws.Range("H119").Value = 31499.715
$ws.Range("J119").Value = 31499.715
$ws.Range("L119").Value = 31499.715
$ws.Range("N119").Value = -41175.715
$ws.Range("H135").Value = 68685.39999999999
$ws.Range("J135").Value = 68685.39999999999
$ws.Range("L135").Value = 68685.39999999999
$ws.Range("N135").Value = -78825.39999999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 33485.098
$ws.Range("I3").Value = 1035.3846
$ws.Range("J3").Value = 56921
$ws.Range("K3").Value = 1035.3846
$ws.Range("L3").Value = 56921
$ws.Range("M3").Value = -921.3846000000001
$ws.Range("N3").Value = -57149
$ws.Range("H35").Value = 39800
$ws.Range("J35").Value = 39800
$ws.Range("L35").Value = 39800
$ws.Range("N35").Value = -40420
$ws.Range("H82").Value = 101428.57
$ws.Range("I82").Value = 101428.57
$ws.Range("K82").Value = 101428.57
$ws.Range("M82").Value = -101045.57
$ws.Range("H85").Value = 101428.57
$ws.Range("I85").Value = 101428.57
$ws.Range("K85").Value = 101428.57
$ws.Range("M85").Value = -100102.57
$ws.Range("H86").Value = 169105.25
$ws.Range("I86").Value = 3222.2856
$ws.Range("J86").Value = 401341.4
$ws.Range("K86").Value = 3222.2856
$ws.Range("L86").Value = 401341.4
$ws.Range("M86").Value = -2099.2856
$ws.Range("N86").Value = -403587.4
$ws.Range("H89").Value = 169105.25
$ws.Range("I89").Value = 3222.2856
$ws.Range("J89").Value = 401341.4
$ws.Range("K89").Value = 16111.428
$ws.Range("L89").Value = 2006707
$ws.Range("M89").Value = -10495.428
$ws.Range("N89").Value = -2017939
$ws.Range("H107").Value = 32671.176
$ws.Range("I107").Value = 41493.08
$ws.Range("K107").Value = 41493.08
$ws.Range("M107").Value = -39573.08
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2531
$ws.Range("I31").Value = 1647.3846
$ws.Range("J31").Value = 3414.6155
$ws.Range("K31").Value = 1647.3846
$ws.Range("L31").Value = 3414.6155
$ws.Range("M31").Value = -1352.3846
$ws.Range("N31").Value = -4004.6155
$ws.Range("H34").Value = 2531
$ws.Range("I34").Value = 1647.3846
$ws.Range("J34").Value = 3414.6155
$ws.Range("K34").Value = 1647.3846
$ws.Range("L34").Value = 3414.6155
$ws.Range("M34").Value = -1445.3846
$ws.Range("N34").Value = -3818.6155
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1388.4025
$ws.Range("I131").Value = 411.25
$ws.Range("J131").Value = 1494.0405
$ws.Range("K131").Value = 1233.75
$ws.Range("L131").Value = 4482.1215
$ws.Range("M131").Value = 3806.25
$ws.Range("N131").Value = -14562.1215
$ws.Range("H134").Value = 3183.375
$ws.Range("I134").Value = 1239.4762
$ws.Range("J134").Value = 5331.8945
$ws.Range("K134").Value = 3718.4286
$ws.Range("L134").Value = 15995.6835
$ws.Range("M134").Value = 1351.5714
$ws.Range("N134").Value = -26135.6835
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2625.3333
$ws.Range("I126").Value = 1400
$ws.Range("J126").Value = 2778.5
$ws.Range("K126").Value = 4200
$ws.Range("L126").Value = 8335.5
$ws.Range("M126").Value = -1730
$ws.Range("N126").Value = -13275.5
$ws.Range("H132").Value = 2495.9592
$ws.Range("I132").Value = 2146.5833
$ws.Range("J132").Value = 3463.4614
$ws.Range("K132").Value = 6439.749899999999
$ws.Range("L132").Value = 10390.3842
$ws.Range("M132").Value = -3909.749899999999
$ws.Range("N132").Value = -15450.3842
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 275307.5
$ws.Range("J36").Value = 275307.5
$ws.Range("L36").Value = 275307.5
$ws.Range("N36").Value = -276431.5
$ws.Range("H132").Value = 3176.75
$ws.Range("I132").Value = 2112.5833
$ws.Range("J132").Value = 4773
$ws.Range("K132").Value = 6337.749899999999
$ws.Range("L132").Value = 14319
$ws.Range("M132").Value = -3807.749899999999
$ws.Range("N132").Value = -19379
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 15000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 15000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 15000
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -16020
$ws.Range("H62").Value = 4600
$ws.Range("I62").Value = 3933.3333
$ws.Range("K62").Value = 3933.3333
$ws.Range("M62").Value = -3309.3333
$ws.Range("H65").Value = 4600
$ws.Range("I65").Value = 3933.3333
$ws.Range("K65").Value = 19666.6665
$ws.Range("M65").Value = -16546.6665
$ws.Range("H113").Value = 400.25
$ws.Range("I113").Value = 300.27274
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 900.81822
$ws.Range("L113").Value = 4500
$ws.Range("M113").Value = 1269.18178
$ws.Range("N113").Value = -8840
$ws.Range("H132").Value = 1890.9131
$ws.Range("I132").Value = 1146.6471
$ws.Range("K132").Value = 3439.9413
$ws.Range("M132").Value = -909.9412999999995
$ws.Range("H136").Value = 1862.1904
$ws.Range("I136").Value = 1755.3
$ws.Range("K136").Value = 5265.9
$ws.Range("M136").Value = -2715.9
